$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new handback file finished processing (a01eb00f-...), and the existing
# in-flight file got its generated id/hash refreshed (a5cfe79f... ->
# 9606c80f...). This adds a second data row to every sheet.
# ---------------------------------------------------------------------------
$newId    = "9606c80f-3135-4c57-8498-146ae6416c79"
$secondId = "a01eb00f-f7c8-4c83-b44f-fa64fe24f78b"
$zhHash   = "93c1ca41ff21c0ce0c8e5e60f097a637f920f49e"
$deHash   = "438042bd86ac7f1270d1a497c98a8187f54661e7"

$repoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a02f2c01a76ca446607a0ef475c3500a7439b9bc/e2e"
$zhRepoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a906eaebb4e10a20dd311d021b6779c282ad65de/e2e"
$deRepoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7fa64c3a29facbcc6eb2a611e07c0a0dcb38e89f/e2e"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Update existing row 2 (first record) with refreshed file name / timestamp.
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-09-07 09:35:24"

# Add new row 3 (second record).
$wsOverview.Range("A3").Value = "$secondId.md"
$wsOverview.Range("B3").Value = "e2e\$secondId.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-07 09:35:24"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild hyperlinks (engine quirk: Hyperlinks.Delete() on any range clears
# the whole sheet's collection, so clear once then re-add everything needed).
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$repoUrl/$newId.md", "", "", "e2e\$newId.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$repoUrl/$secondId.md", "", "", "e2e\$secondId.md") | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# Update existing row 2 (first record).
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "$newId.$zhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-07 09:35:15"
$wsZh.Range("I2").Value = "$newId.md"
$wsZh.Range("J2").Value = "$newId.$zhHash.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-07 09:35:44"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""

# Add new row 3 (second record).
$wsZh.Range("A3").Value = "$secondId.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$secondId.$deHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-07 09:35:15"
$wsZh.Range("I3").Value = "$secondId.md"
$wsZh.Range("J3").Value = "$secondId.$deHash.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-07 09:35:44"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoUrl/$newId.md", "", "", "$newId.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$zhRepoUrl/$newId.md", "", "", "$newId.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoUrl/$secondId.md", "", "", "$secondId.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$zhRepoUrl/$secondId.md", "", "", "$secondId.md") | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

# Update existing row 2 (first record).
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "$newId.$zhHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-07 09:35:24"
$wsDe.Range("I2").Value = "$newId.md"
$wsDe.Range("J2").Value = "$newId.$zhHash.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-07 09:35:53"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""

# Add new row 3 (second record).
$wsDe.Range("A3").Value = "$secondId.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$secondId.$deHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-07 09:35:24"
$wsDe.Range("I3").Value = "$secondId.md"
$wsDe.Range("J3").Value = "$secondId.$deHash.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-07 09:35:53"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoUrl/$newId.md", "", "", "$newId.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$deRepoUrl/$newId.md", "", "", "$newId.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoUrl/$secondId.md", "", "", "$secondId.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$deRepoUrl/$secondId.md", "", "", "$secondId.md") | Out-Null
